# Append two new log rows (rows 20 and 21) to Sheet1, matching the
# existing "run_id, rss_url_id, date, response, item_count" data layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 20
$ws.Cells.Item(20, 1).Value = 19
$ws.Cells.Item(20, 2).Value = 1
$ws.Cells.Item(20, 3).Value = "2024-06-15 07:12:54"
$ws.Cells.Item(20, 4).Value = 200
$ws.Cells.Item(20, 5).Value = 3

# Row 21
$ws.Cells.Item(21, 1).Value = 20
$ws.Cells.Item(21, 2).Value = 2
$ws.Cells.Item(21, 3).Value = "2024-06-15 07:12:54"
$ws.Cells.Item(21, 4).Value = 200
$ws.Cells.Item(21, 5).Value = 0
